# Generate Report for Handoff
#
# Replaces the old source files (78b55644-... / 7b2c538e-...) with the new
# handoff batch (230d5acb-... / ffffd46d8d01-...), flips the status from
# "Handed back: in sync with en-US" to "Ready for handoff", clears the
# now-stale "Latest Target File" / "Latest Handback File" / hyperlinks, and
# points zh-cn + de-de at the freshly generated xliff for the new source.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.md"
$ov.Range("B2").Value = "e2e\230d5acb-504b-4637-9f8e-a5c285449e76.md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-09-04 07:09:09"

$ov.Range("A3").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$ov.Range("B3").Value = "e2e\ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-04 07:09:09"

# Columns shrink now that the status text is shorter.
$ov.Columns.Item(5).ColumnWidth = 16.333333333333332
$ov.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-04 07:09:00"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-04 07:09:00"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

# The handback hyperlinks (I2 / I3) are gone now that the file is merely
# queued for handoff, not handed back. This host only supports wiping the
# whole collection at once, so drop everything and re-add just the two
# source-file links (A2 / A3) that should still exist, preserving their
# original targets and display text.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/230d5acb-504b-4637-9f8e-a5c285449e76.md", [Type]::Missing, [Type]::Missing, "230d5acb-504b-4637-9f8e-a5c285449e76.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", [Type]::Missing, [Type]::Missing, "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md")

$zh.Columns.Item(3).ColumnWidth = 16.333333333333332
$zh.Columns.Item(9).ColumnWidth = 17.833333333333332
$zh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf"
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""

$de.Range("A3").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf"
$de.Range("K3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/230d5acb-504b-4637-9f8e-a5c285449e76.md", [Type]::Missing, [Type]::Missing, "230d5acb-504b-4637-9f8e-a5c285449e76.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", [Type]::Missing, [Type]::Missing, "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md")

$de.Columns.Item(3).ColumnWidth = 16.333333333333332
$de.Columns.Item(9).ColumnWidth = 17.833333333333332
$de.Columns.Item(10).ColumnWidth = 20.833333333333332
